# "Modificaciones en reportes y en lista de tareas"
#
# - Adds a currency/percent-style mark (reusing existing number format) on C55
#   (task "Logueo de aplicación (configurable)" / Lucas) to flag it, matching
#   the "$ en los totales" follow-up task added below.
# - Adds 3 new task rows to the "lista de tareas" sheet (Hoja1):
#     58: Agregar los comandos abajo, en la pantalla (teclas rapidas) - Lucas
#     59: Agregar signo $ en los totales            (replaces old row 59 text)
#     60: En reportes mensual y anual poner mes y año, no dia! - Agustina
# - The task that used to live in row 59 ("Ivan: preguntar reportes...") is
#   pushed down to row 66.
# - Updates the window scroll position / selection to the new working area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Mark C55 with the same percent-style formatting used by the other "done"
# flags in column C (reuses existing style, cell itself stays empty).
$ws.Range("C55").NumberFormat = "0%"

# New task row 58
$ws.Range("A58").Value = "Agregar los comandos abajo, en la pantalla (teclas rapidas)"
$ws.Range("B58").Value = "Lucas"

# Row 59 now holds the new "$ en los totales" task (overwrites previous text)
$ws.Range("A59").Value = "Agregar signo `$ en los totales"

# New task row 60
$ws.Range("A60").Value = "En reportes mensual y anual poner mes y año, no dia!"
$ws.Range("B60").Value = "Agustina"

# Former row 59 content is relocated further down to row 66
$ws.Range("A66").Value = "Ivan: preguntar reportes - preguntar autorizacion requerida en que funciones - preguntar login"

# Reflect the new scroll position / active selection used while editing
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
$ws.Range("B59").Select()
